$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Name column (A2:A10) with new landmark names
$ws.Range("A2").Value = "Bellagio Hotel"
$ws.Range("A3").Value = "The Getty"
$ws.Range("A4").Value = "Flatiron"
$ws.Range("A5").Value = "KOIN Center"
$ws.Range("A6").Value = "The Parthenon"
$ws.Range("A7").Value = "Olympia Theater"
$ws.Range("A8").Value = "Space Needle"
$ws.Range("A9").Value = "Coors Field"
$ws.Range("A10").Value = "Anson Mills"

# Re-fit column A width to the new (longer) content
$ws.Columns.Item(1).ColumnWidth = 15

# Update the selected cell shown when the workbook is opened
$ws.Range("L7").Select() | Out-Null

# Restore window position metadata
$excel.ActiveWindow.Left = 3900
$excel.ActiveWindow.Top = 3900
